$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking values (e.g. "1.001", "328.00")
# are preserved exactly as text, matching the source data which stores them as strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.500.84"
$ws.Range("E2").Value = "  +1.55%  "

$ws.Range("D3").Value = "1.908.91"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").Value = "328.00"
$ws.Range("E5").Value = "  -1.54%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").Value = "0.4659"
$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("D8").Value = "0.4081"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "47.74"
$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").Value = "0.08003"
$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").Value = "1.006"
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").Value = "22.29"
$ws.Range("E12").Value = "  +2.54%  "

$ws.Range("D13").Value = "1.906.61"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").Value = "5.942"
$ws.Range("E14").Value = "  -0.14%  "

$ws.Range("D15").Value = "7.123"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").Value = "89.04"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.06595"
$ws.Range("E18").Value = "  +0.77%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.00001028"
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("D20").Value = "17.72"
$ws.Range("E20").Value = "  +1.27%  "

$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "29.512.31"
$ws.Range("E22").Value = "  +1.71%  "

$ws.Range("D23").Value = "5.533"
$ws.Range("E23").Value = "  +1.33%  "

$ws.Range("D24").Value = "11.49"
$ws.Range("E24").Value = "  +2.91%  "

$ws.Range("D25").Value = "2.206"
$ws.Range("E25").Value = "  -1.66%  "

$ws.Range("D26").Value = "2.128.48"
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").Value = "153.61"
$ws.Range("E27").Value = "  -2.79%  "

$ws.Range("D28").Value = "19.76"
$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").Value = "5.726"
$ws.Range("E29").Value = "  +6.28%  "

$ws.Range("D30").Value = "2.122"
$ws.Range("E30").Value = "  +0.75%  "

$ws.Range("D31").Value = "116.68"
$ws.Range("E31").Value = "  -2.09%  "

$ws.Range("D32").Value = "1.074"
$ws.Range("E32").Value = "  +9.26%  "

$ws.Range("D33").Value = "0.09453"
$ws.Range("E33").Value = "  +0.72%  "

$ws.Range("D34").Value = "1.422"
$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("D35").Value = "3.575"
$ws.Range("E35").Value = "  -0.69%  "

$ws.Range("D36").Value = "5.377"
$ws.Range("E36").Value = "  +1.29%  "

$ws.Range("D37").Value = "0.02256"
$ws.Range("E37").Value = "  +1.22%  "

$ws.Range("D38").Value = "0.06078"
$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("D39").Value = "8.380"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").Value = "1.172"
$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("D41").Value = "0.5863"
$ws.Range("E41").Value = "  +0.68%  "

$ws.Range("D42").Value = "0.1832"
$ws.Range("E42").Value = "  +0.31%  "

$ws.Range("D43").Value = "10.09"
$ws.Range("E43").Value = "  -0.95%  "

$ws.Range("D44").Value = "1.302"
$ws.Range("E44").Value = "  +3.41%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.07737"
$ws.Range("E45").Value = "  +9.97%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "2.372"
$ws.Range("E46").Value = "  +1.37%  "

$ws.Range("D47").Value = "12.16"
$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("D48").Value = "0.5537"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("D49").Value = "1.922"
$ws.Range("E49").Value = "  +0.84%  "

$ws.Range("D50").Value = "113.26"
$ws.Range("E50").Value = "  +0.69%  "

$ws.Range("D51").Value = "0.2928"
$ws.Range("E51").Value = "  +4.64%  "
